# "bo sung qd can cu" - split the "Căn cứ Luật Quản lý thuế ngày 13 tháng 6
# năm 2019;" sentence so the fixed date is replaced by a <luat_qlt_ngay>
# placeholder token, and relocate the document's (hidden) _GoBack bookmark
# from the very end of the document to right after that new placeholder
# run, matching the target OOXML exactly.

$d = $word.ActiveDocument

# The document currently carries a single _GoBack bookmark at the very end
# (right after the last sentence of the document). It needs to move to the
# newly split run below, so drop it here; we recreate it (with a fresh id)
# in the InsertXML payload used below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Locate the whole old sentence so we can replace it wholesale with the
# three-run + bookmark structure the diff wants.
$found = $d.Content.Duplicate
$found.Find.Execute(
    "Căn cứ Luật Quản lý thuế ngày 13 tháng 6 năm 2019;",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$target = $d.Range($found.Start, $found.End)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
            '<pkg:xmlData>' +
                '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                    '<w:body>' +
                        '<w:p>' +
                            '<w:r>' +
                                '<w:rPr><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="nl-NL"/></w:rPr>' +
                                '<w:t xml:space="preserve">Căn cứ Luật Quản lý thuế </w:t>' +
                            '</w:r>' +
                            '<w:r>' +
                                '<w:rPr><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="nl-NL"/></w:rPr>' +
                                '<w:t>&lt;luat_qlt_ngay&gt;</w:t>' +
                            '</w:r>' +
                            '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
                            '<w:bookmarkEnd w:id="0"/>' +
                            '<w:r>' +
                                '<w:rPr><w:i/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="nl-NL"/></w:rPr>' +
                                '<w:t>;</w:t>' +
                            '</w:r>' +
                        '</w:p>' +
                    '</w:body>' +
                '</w:document>' +
            '</pkg:xmlData>' +
        '</pkg:part>' +
    '</pkg:package>'

$target.InsertXML($xml)
